$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44922; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    3  = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    4  = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    5  = @{ D = 44895; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    6  = @{ D = 44943; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    7  = @{ D = 44959; J = 30; K = 19000; L = 19000; M = 19000; P = 1462 }
    8  = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    9  = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    10 = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    11 = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    12 = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    13 = @{ D = 44930; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    14 = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
    15 = @{ D = 44915; J = 50; K = 18000; L = 18000; M = 18000; P = 1385 }
    16 = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    17 = @{ D = 44894; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    18 = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    19 = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
}

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $epoch.AddDays([double]$vals.D)
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}

$wb.Save()
